$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data appended after the last existing row (row 79 -> row 80).
$row = 80

# Column A holds a date formatted as text (e.g. "09/03/2025"), stored as a
# plain string, not an Excel date serial. Writing the literal string via
# .Value would make Excel auto-recognize it as a date and convert it to a
# serial number, which both changes the cell type and allocates a new
# NumberFormat style. To keep the cell a genuine text value (matching every
# other row in the column, with no style override), build the text in a
# scratch cell using TEXT(), copy it, and paste *values only* into place.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=TEXT(DATE(2025,11,20),""mm/dd/yyyy"")"
$scratch.Copy()
$ws.Range("A80").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Cells.Item($row, 2).Value = 0.2044597188475565
$ws.Cells.Item($row, 3).Value = 0.7955402811524435
